# Weekly update: a new price record (week of 2022-01-25) is added at the
# top of the "Macroferia Regional de Talca - Ciruela" data block (which is
# sorted most-recent-first starting at row 37). All the existing records
# from row 37 down to row 73 shift down by one row (to rows 38-74), and the
# new record is written into the now-empty row 37.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 37; Excel shifts rows 37:73 down to 38:74 and carries
# the row formatting (e.g. the date style on column D) along with them.
$ws.Rows("37:37").Insert()

# Populate the newly-opened row 37 with the new weekly record.
$ws.Range("A37").Value = 5
$ws.Range("B37").Value = "Macroferia Regional de Talca"
$ws.Range("C37").Value = "Maule"
$ws.Range("D37").Value = 44586
$ws.Range("E37").Value = 7
$ws.Range("F37").Value = "Fruta"
$ws.Range("G37").Value = 100103
$ws.Range("H37").Value = "Frutos de hueso (carozo)"
$ws.Range("I37").Value = 100103002
$ws.Range("J37").Value = "Ciruela"
$ws.Range("K37").Value = "Black Amber"
$ws.Range("L37").Value = "Primera"
$ws.Range("M37").Value = 300
$ws.Range("N37").Value = 9000
$ws.Range("O37").Value = 9000
$ws.Range("P37").Value = 9000
$ws.Range("Q37").Value = "$/bandeja 18 kilos granel"
$ws.Range("R37").Value = "Provincia de Curicó"
$ws.Range("S37").Value = 500
$ws.Range("T37").Value = 18
